$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.739.18'
$ws.Range('E2').Value = '  -1.17%  '

$ws.Range('D3').Value = '1.627.19'
$ws.Range('E3').Value = '  -1.07%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.41'
$ws.Range('E5').Value = '  +0.24%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5106'
$ws.Range('E6').Value = '  +0.10%  '

$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2561'
$ws.Range('E8').Value = '  -0.09%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06318'
$ws.Range('E9').Value = '  -0.64%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.41'
$ws.Range('E10').Value = '  -0.80%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07778'
$ws.Range('E11').Value = '  +0.17%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.234'
$ws.Range('E12').Value = '  -1.40%  '

$ws.Range('D13').Value = '1.630.74'
$ws.Range('E13').Value = '  -0.83%  '

$ws.Range('D14').Value = '1.850.59'
$ws.Range('E14').Value = '  -1.23%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5503'
$ws.Range('E15').Value = '  +1.18%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.57'
$ws.Range('E16').Value = '  -1.14%  '

$ws.Range('D17').Value = '0.0₅7537'
$ws.Range('E17').Value = '  -2.33%  '

$ws.Range('D18').Value = '25.792.33'
$ws.Range('E18').Value = '  -1.05%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.002'

$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.402'
$ws.Range('E20').Value = '  -0.66%  '

$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '193.91'
$ws.Range('E21').Value = '  -2.52%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.824'
$ws.Range('E22').Value = '  -1.16%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.995'
$ws.Range('E23').Value = '  -0.86%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.004'
$ws.Range('E24').Value = '  +0.03%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.888'
$ws.Range('E25').Value = '  +0.89%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.06'
$ws.Range('E26').Value = '  +0.40%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.54'
$ws.Range('E28').Value = '  -0.61%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.739'
$ws.Range('E29').Value = '  -1.09%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.238'
$ws.Range('E30').Value = '  +0.19%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04874'
$ws.Range('E31').Value = '  +0.30%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.226'
$ws.Range('E32').Value = '  -0.98%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.169'
$ws.Range('E33').Value = '  +0.10%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.539'
$ws.Range('E34').Value = '  +0.86%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.374'
$ws.Range('E35').Value = '  +0.25%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.8926'
$ws.Range('E36').Value = '  -0.78%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5498'
$ws.Range('E37').Value = '  +0.63%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.535'
$ws.Range('E38').Value = '  -1.97%  '

$ws.Range('D39').Value = '1.112.83'
$ws.Range('E39').Value = '  -2.62%  '

$ws.Range('E40').Value = '  -0.96%  '

$ws.Range('E41').Value = '  -0.06%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.564'
$ws.Range('E42').Value = '  +3.18%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7942'
$ws.Range('E43').Value = '  -2.25%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.16'
$ws.Range('E44').Value = '  -2.14%  '

$ws.Range('D45').Value = '1.776.81'
$ws.Range('E45').Value = '  -0.37%  '

$ws.Range('D46').Value = '0.0₈112'
$ws.Range('E46').Value = '  -12.81%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4430'
$ws.Range('E47').Value = '  -2.24%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9993'
$ws.Range('E48').Value = '  -0.17%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.57'
$ws.Range('E49').Value = '  -0.79%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05133'
$ws.Range('E50').Value = '  +1.51%  '

$ws.Range('E51').Value = '  +2.83%  '
